$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 9972.75
$ws.Range("I51").Value = 8949.5
$ws.Range("K51").Value = 8949.5
$ws.Range("M51").Value = -8465.5
$ws.Range("H69").Value = 2060
$ws.Range("J69").Value = 2000
$ws.Range("L69").Value = 6000
$ws.Range("N69").Value = -7748
$ws.Range("H72").Value = 2060
$ws.Range("J72").Value = 2000
$ws.Range("L72").Value = 18000
$ws.Range("N72").Value = -26736
$ws.Range("H116").Value = 4175.5625
$ws.Range("I116").Value = 1999.4
$ws.Range("K116").Value = 1999.4
$ws.Range("M116").Value = 1442.6

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2047.8
$ws.Range("I2").Value = 1336.8182
$ws.Range("J2").Value = 4003
$ws.Range("K2").Value = 1336.8182
$ws.Range("L2").Value = 4003
$ws.Range("M2").Value = -1223.8182
$ws.Range("N2").Value = -4229
$ws.Range("H74").Value = 2625.1538
$ws.Range("I74").Value = 2071.3
$ws.Range("K74").Value = 2071.3
$ws.Range("M74").Value = -1197.3
$ws.Range("H77").Value = 2625.1538
$ws.Range("I77").Value = 2071.3
$ws.Range("K77").Value = 10356.5
$ws.Range("M77").Value = -5988.5
$ws.Range("H116").Value = 2047.8
$ws.Range("I116").Value = 1336.8182
$ws.Range("J116").Value = 4003
$ws.Range("K116").Value = 1336.8182
$ws.Range("L116").Value = 4003
$ws.Range("M116").Value = 957.1818000000001
$ws.Range("N116").Value = -8591

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2047.8
$ws.Range("I3").Value = 1336.8182
$ws.Range("J3").Value = 4003
$ws.Range("K3").Value = 1336.8182
$ws.Range("L3").Value = 4003
$ws.Range("M3").Value = -1222.8182
$ws.Range("N3").Value = -4231
$ws.Range("H80").Value = 872.6
$ws.Range("I80").Value = 886.25
$ws.Range("J80").Value = 863.5
$ws.Range("K80").Value = 886.25
$ws.Range("L80").Value = 863.5
$ws.Range("M80").Value = 111.75
$ws.Range("N80").Value = -2859.5
$ws.Range("H83").Value = 872.6
$ws.Range("I83").Value = 886.25
$ws.Range("J83").Value = 863.5
$ws.Range("K83").Value = 4431.25
$ws.Range("L83").Value = 4317.5
$ws.Range("M83").Value = 560.75
$ws.Range("N83").Value = -14301.5
$ws.Range("H94").Value = 2271.4211
$ws.Range("I94").Value = 1626.6923
$ws.Range("K94").Value = 1626.6923
$ws.Range("M94").Value = -1175.6923
$ws.Range("H105").Value = 3573778.5
$ws.Range("I105").Value = 2825
$ws.Range("J105").Value = 5002160
$ws.Range("K105").Value = 2825
$ws.Range("L105").Value = 5002160
$ws.Range("M105").Value = -1078
$ws.Range("N105").Value = -5005654

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1183.7222
$ws.Range("I16").Value = 950.7
$ws.Range("J16").Value = 1475
$ws.Range("K16").Value = 950.7
$ws.Range("L16").Value = 1475
$ws.Range("M16").Value = -663.7
$ws.Range("N16").Value = -2049
$ws.Range("H31").Value = 15544.483
$ws.Range("J31").Value = 5034.1177
$ws.Range("L31").Value = 5034.1177
$ws.Range("N31").Value = -5624.1177
$ws.Range("H34").Value = 15544.483
$ws.Range("J34").Value = 5034.1177
$ws.Range("L34").Value = 5034.1177
$ws.Range("N34").Value = -5438.1177
$ws.Range("H88").Value = 37671.5
$ws.Range("J88").Value = 37671.5
$ws.Range("L88").Value = 37671.5
$ws.Range("N88").Value = -38483.5
$ws.Range("H91").Value = 37671.5
$ws.Range("J91").Value = 37671.5
$ws.Range("L91").Value = 37671.5
$ws.Range("N91").Value = -40479.5
$ws.Range("H113").Value = 1183.7222
$ws.Range("I113").Value = 950.7
$ws.Range("J113").Value = 1475
$ws.Range("K113").Value = 950.7
$ws.Range("L113").Value = 1475
$ws.Range("M113").Value = 1219.3
$ws.Range("N113").Value = -5815
$ws.Range("H134").Value = 849.3333
$ws.Range("I134").Value = 724.2857
$ws.Range("J134").Value = 1287
$ws.Range("K134").Value = 2172.8571
$ws.Range("L134").Value = 3861
$ws.Range("M134").Value = 362.1428999999998
$ws.Range("N134").Value = -8931

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 543
$ws.Range("I98").Value = 322.5
$ws.Range("K98").Value = 967.5
$ws.Range("M98").Value = 530.5
$ws.Range("H131").Value = 750.97
$ws.Range("J131").Value = 750.97
$ws.Range("L131").Value = 2252.91
$ws.Range("N131").Value = -12332.91

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 465.57144
$ws.Range("J107").Value = 622.25
$ws.Range("L107").Value = 622.25
$ws.Range("N107").Value = -4462.25
$ws.Range("H113").Value = 2638.2778
$ws.Range("I113").Value = 1891.909
$ws.Range("J113").Value = 3811.1428
$ws.Range("K113").Value = 1891.909
$ws.Range("L113").Value = 3811.1428
$ws.Range("M113").Value = 278.0909999999999
$ws.Range("N113").Value = -8151.1428
$ws.Range("H122").Value = 1411.4546
$ws.Range("I122").Value = 1280.6666
$ws.Range("K122").Value = 3841.9998
$ws.Range("M122").Value = -1391.9998

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 576266.4
$ws.Range("I132").Value = 805512.75
$ws.Range("J132").Value = 3150.5
$ws.Range("K132").Value = 2416538.25
$ws.Range("L132").Value = 9451.5
$ws.Range("M132").Value = -2414008.25
$ws.Range("N132").Value = -14511.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 1000004
$ws.Range("I17").Value = 1000004
$ws.Range("K17").Value = 1000004
$ws.Range("M17").Value = -999832
$ws.Range("H81").Value = 1857.1428
$ws.Range("I81").Value = 1857.1428
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 3714.2856
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -2653.2856
$ws.Range("H84").Value = 1857.1428
$ws.Range("I84").Value = 1857.1428
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 18571.428
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -13267.428
$ws.Range("H110").Value = 25000
$ws.Range("J110").Value = 25000
$ws.Range("L110").Value = 25000
$ws.Range("N110").Value = -33180
$ws.Range("H122").Value = 2018.7693
$ws.Range("I122").Value = 1790
$ws.Range("J122").Value = 2781.3333
$ws.Range("K122").Value = 5370
$ws.Range("L122").Value = 8343.999899999999
$ws.Range("M122").Value = -2920
$ws.Range("N122").Value = -13243.9999
$ws.Range("H126").Value = 747.9583
$ws.Range("I126").Value = 660.55554
$ws.Range("J126").Value = 1010.1667
$ws.Range("K126").Value = 1981.66662
$ws.Range("L126").Value = 3030.5001
$ws.Range("M126").Value = 488.33338
$ws.Range("N126").Value = -7970.5001

# --- Cell removals (ClearContents) ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("N81").ClearContents()
$ws.Range("N84").ClearContents()
